$d = $word.ActiveDocument

# Locate the paragraph that ends with "...terraform first and trying to
# combine deployments instead of one resource at a time." — it is the last
# paragraph of real content before the trailing run of blank paragraphs.
$anchorText = "terraform first and trying to combine deployments instead of one resource at a time."
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "*$anchorText*") {
        $anchorIndex = $i
        $found = $true
    }
}

$anchorPara = $d.Paragraphs($anchorIndex)
$r = $anchorPara.Range
$r.Collapse(0)

# Build the new block: a blank spacer paragraph, the centered/bold/underlined
# section heading, and the journal-entry paragraph (the four sentences are
# joined into one insertion since they share identical run formatting).
$heading = "AWS practice/Linux refamiliarization / SSH configuration"
$bodyPart1 = "Today I downloaded and installed openssh using ps command line and created ssh key pairs in aws gui for 2 instances "
$bodyPart2 = "using RHEL. I was able to mess around with port 22 and get into both systems and perform  sudo yum update commands. I also played around with directories but I think lack of access to port 80/443 when I setup the instances is what is preventing "
$bodyPart3 = "things like ansible, ufw, and nano from being installed."
$bodyPart4 = " Also went over some netapp storage stuff."
$body = $bodyPart1 + $bodyPart2 + $bodyPart3 + $bodyPart4

$newBlock = "`r" + $heading + "`r" + $body
$r.InsertAfter($newBlock)

# Re-resolve the three freshly-inserted paragraphs by index.
$headingPara = $d.Paragraphs($anchorIndex + 2)
$bodyPara = $d.Paragraphs($anchorIndex + 3)

# Center the heading paragraph and make it bold/underlined, matching the
# "AWS practice/Linux refamiliarization / SSH configuration" title style
# used elsewhere in the document.
$headingPara.Alignment = 1
$hr = $headingPara.Range
$hr.Font.Bold = 1
$hr.Font.BoldBi = 1
$hr.Font.Underline = 1

Write-Output "Inserted AWS practice/Linux refamiliarization / SSH configuration section after paragraph $anchorIndex"
